$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I, J), styled like the other headers (same style as H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for rows 2-62: Row number, I value, J value
$data = @(
    @{Row=2; I=7; J=7},
    @{Row=3; I=2; J=4},
    @{Row=4; I=10; J=10},
    @{Row=5; I=7; J=7},
    @{Row=6; I=4; J=5},
    @{Row=7; I=7; J=8},
    @{Row=8; I=6; J=7},
    @{Row=9; I=8; J=8},
    @{Row=10; I=7; J=7},
    @{Row=11; I=7; J=7},
    @{Row=12; I=7; J=7},
    @{Row=13; I=7; J=7},
    @{Row=14; I=7; J=7},
    @{Row=15; I=8; J=8},
    @{Row=16; I=5; J=6},
    @{Row=17; I=8; J=8},
    @{Row=18; I=8; J=8},
    @{Row=19; I=6; J=6},
    @{Row=20; I=6; J=7},
    @{Row=21; I=7; J=7},
    @{Row=22; I=8; J=8},
    @{Row=23; I=6; J=6},
    @{Row=24; I=10; J=10},
    @{Row=25; I=7; J=7},
    @{Row=26; I=8; J=8},
    @{Row=27; I=5; J=6},
    @{Row=28; I=1; J=4},
    @{Row=29; I=5; J=5},
    @{Row=30; I=6; J=7},
    @{Row=31; I=8; J=8},
    @{Row=32; I=8; J=9},
    @{Row=33; I=6; J=6},
    @{Row=34; I=5; J=5},
    @{Row=35; I=7; J=8},
    @{Row=36; I=10; J=10},
    @{Row=37; I=7; J=8},
    @{Row=38; I=11; J=11},
    @{Row=39; I=9; J=9},
    @{Row=40; I=8; J=8},
    @{Row=41; I=9; J=9},
    @{Row=42; I=8; J=8},
    @{Row=43; I=9; J=9},
    @{Row=44; I=8; J=9},
    @{Row=45; I=7; J=7},
    @{Row=46; I=7; J=7},
    @{Row=47; I=7; J=8},
    @{Row=48; I=7; J=8},
    @{Row=49; I=5; J=5},
    @{Row=50; I=6; J=7},
    @{Row=51; I=8; J=8},
    @{Row=52; I=8; J=8},
    @{Row=53; I=11; J=11},
    @{Row=54; I=7; J=8},
    @{Row=55; I=7; J=7},
    @{Row=56; I=4; J=5},
    @{Row=57; I=7; J=8},
    @{Row=58; I=7; J=8},
    @{Row=59; I=9; J=9},
    @{Row=60; I=4; J=4},
    @{Row=61; I=5; J=5},
    @{Row=62; I=5; J=5}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 9).Value = $entry.I
    $ws.Cells.Item($entry.Row, 10).Value = $entry.J
}
